$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.561.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.54%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.988.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.16%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '381.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.03%  '

$ws.Range('E7').Value = '  +2.12%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.33%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0860'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.461.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.29%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.11%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.76%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.985.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.60%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.995'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.585.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.62%  '

$ws.Range('E20').Value = '  +1.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.67%  '

$ws.Range('E22').Value = '  +1.43%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.31%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.94%  '

$ws.Range('E25').Value = '  +2.72%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.20%  '

$ws.Range('E27').Value = '  +5.00%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.51%  '

$ws.Range('E31').Value = '  +0.95%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.57%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.22%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0445'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.76%  '

$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.76%  '

$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('E38').Value = '  +6.31%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.93'
$ws.Range('D39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.20%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '129.09'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.19%  '

$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.116'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.70%  '

$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +14.29%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.62%  '

$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.269'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.47%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.35%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.031.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.285.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.94%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.534'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +18.89%  '
